$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws "D2" "71.853.42"
Set-TextValue $ws "E2" "  -0.76%  "
Set-TextValue $ws "D3" "2.648.45"
Set-TextValue $ws "E3" "  +0.99%  "
Set-TextValue $ws "E4" "  +0.04%  "
Set-TextValue $ws "D5" "597.00"
Set-TextValue $ws "E5" "  -1.09%  "
Set-TextValue $ws "D6" "174.41"
Set-TextValue $ws "E6" "  -2.32%  "
Set-TextValue $ws "E7" "  +0.04%  "
Set-TextValue $ws "D8" "0.522"
Set-TextValue $ws "E8" "  -0.59%  "
Set-TextValue $ws "D9" "2.643.82"
Set-TextValue $ws "E9" "  +0.91%  "
Set-TextValue $ws "E10" "  -1.52%  "
Set-TextValue $ws "E11" "  +2.36%  "
Set-TextValue $ws "D12" "0.355"
Set-TextValue $ws "E12" "  +0.89%  "
Set-TextValue $ws "E13" "  -0.87%  "
Set-TextValue $ws "D14" "3.149.17"
Set-TextValue $ws "E14" "  +0.34%  "
Set-TextValue $ws "E15" "  -2.03%  "
Set-TextValue $ws "D16" "71.668.36"
Set-TextValue $ws "E16" "  -0.83%  "
Set-TextValue $ws "D17" "26.19"
Set-TextValue $ws "E17" "  -1.30%  "
Set-TextValue $ws "D18" "2.651.35"
Set-TextValue $ws "E18" "  +1.02%  "
Set-TextValue $ws "D19" "12.17"
Set-TextValue $ws "E19" "  +5.28%  "
Set-TextValue $ws "E20" "  +1.99%  "
Set-TextValue $ws "D21" "369.99"
Set-TextValue $ws "E21" "  -3.33%  "
Set-TextValue $ws "E22" "  -0.22%  "
Set-TextValue $ws "D23" "2.03"
Set-TextValue $ws "E23" "  +1.52%  "
Set-TextValue $ws "D24" "71.96"
Set-TextValue $ws "E24" "  -1.20%  "
Set-TextValue $ws "E25" "  +0.02%  "
Set-TextValue $ws "E26" "  -1.32%  "
Set-TextValue $ws "E27" "  -1.13%  "
Set-TextValue $ws "D28" "2.785.48"
Set-TextValue $ws "E28" "  +1.16%  "
Set-TextValue $ws "E29" "  +0.16%  "
Set-TextValue $ws "D30" "0.0₃0966"
Set-TextValue $ws "E30" "  +1.38%  "
Set-TextValue $ws "E31" "  +0.12%  "
Set-TextValue $ws "D32" "500.34"
Set-TextValue $ws "E32" "  -4.70%  "
Set-TextValue $ws "E33" "  -2.62%  "
Set-TextValue $ws "E34" "  -0.39%  "
Set-TextValue $ws "D35" "0.998"
Set-TextValue $ws "E35" "  -0.09%  "
Set-TextValue $ws "D36" "162.93"
Set-TextValue $ws "E36" "  -0.77%  "
Set-TextValue $ws "D37" "19.48"
Set-TextValue $ws "E37" "  +0.98%  "
Set-TextValue $ws "D38" "18.99"
Set-TextValue $ws "E38" "  -0.61%  "
Set-TextValue $ws "E39" "  -1.11%  "
Set-TextValue $ws "E40" "  -1.67%  "
Set-TextValue $ws "D41" "1.77"
Set-TextValue $ws "E41" "  -3.10%  "
Set-TextValue $ws "E42" "  +0.03%  "
Set-TextValue $ws "E43" "  -1.28%  "
Set-TextValue $ws "E44" "  -0.99%  "
Set-TextValue $ws "E45" "  -0.09%  "

# Row 46: Aave -> OKB
Set-TextValue $ws "B46" "OKB"
Set-TextValue $ws "C46" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws "D46" "39.46"
Set-TextValue $ws "E46" "  -0.05%  "

# Row 47: OKB -> Aave
Set-TextValue $ws "B47" "Aave"
Set-TextValue $ws "C47" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws "D47" "155.68"
Set-TextValue $ws "E47" "  +3.64%  "

Set-TextValue $ws "D48" "0.558"
Set-TextValue $ws "E48" "  +2.89%  "
Set-TextValue $ws "E49" "  +1.04%  "
Set-TextValue $ws "D50" "1.72"
Set-TextValue $ws "E50" "  +1.71%  "

# Row 51: Mantle -> Cronos
Set-TextValue $ws "B51" "Cronos"
Set-TextValue $ws "C51" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws "D51" "0.0752"
Set-TextValue $ws "E51" "  -1.63%  "
